$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3275593.8
$ws.Range("I6").Value = 6174010.5
$ws.Range("J6").Value = 14874.875
$ws.Range("K6").Value = 18522031.5
$ws.Range("L6").Value = 44624.625
$ws.Range("M6").Value = -18521919.5
$ws.Range("N6").Value = -44848.625

$ws.Range("H112").Value = 3856.25
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3856.25
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 11568.75
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -13784.75

$ws.Range("H123").Value = 35920
$ws.Range("J123").Value = 35920
$ws.Range("L123").Value = 35920
$ws.Range("N123").Value = -45720

$ws.Range("H129").Value = 994.23334
$ws.Range("I129").Value = 1293.5385
$ws.Range("J129").Value = 943.7012999999999
$ws.Range("K129").Value = 3880.6155
$ws.Range("L129").Value = 2831.1039
$ws.Range("M129").Value = 1119.3845
$ws.Range("N129").Value = -12831.1039

$ws.Range("H137").Value = 2963837.5
$ws.Range("I137").Value = 9616365
$ws.Range("J137").Value = 7158.4443
$ws.Range("K137").Value = 28849095
$ws.Range("L137").Value = 21475.3329
$ws.Range("M137").Value = -28846545
$ws.Range("N137").Value = -26575.3329

$ws.Range("H141").Value = 3164.5881
$ws.Range("I141").Value = 2228.5833
$ws.Range("J141").Value = 5411
$ws.Range("K141").Value = 6685.749899999999
$ws.Range("L141").Value = 16233
$ws.Range("M141").Value = -1505.749899999999
$ws.Range("N141").Value = -26593

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2692.7273
$ws.Range("I61").Value = 1991.3158
$ws.Range("K61").Value = 1991.3158
$ws.Range("M61").Value = -1779.3158

$ws.Range("H136").Value = 2692.7273
$ws.Range("I136").Value = 1991.3158
$ws.Range("K136").Value = 5973.9474
$ws.Range("M136").Value = -3423.9474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3582.5
$ws.Range("I105").Value = 2713.1667
$ws.Range("K105").Value = 2713.1667
$ws.Range("M105").Value = -966.1667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 929.8889
$ws.Range("I16").Value = 920.3125
$ws.Range("J16").Value = 1006.5
$ws.Range("K16").Value = 920.3125
$ws.Range("L16").Value = 1006.5
$ws.Range("M16").Value = -633.3125
$ws.Range("N16").Value = -1580.5

$ws.Range("H22").Value = 1155.8334
$ws.Range("I22").Value = 269.6154
$ws.Range("J22").Value = 2203.182
$ws.Range("K22").Value = 269.6154
$ws.Range("L22").Value = 2203.182
$ws.Range("M22").Value = 80.38459999999998
$ws.Range("N22").Value = -2903.182

$ws.Range("H31").Value = 4333191
$ws.Range("I31").Value = 1627.3846
$ws.Range("J31").Value = 8778743
$ws.Range("K31").Value = 1627.3846
$ws.Range("L31").Value = 8778743
$ws.Range("M31").Value = -1332.3846
$ws.Range("N31").Value = -8779333

$ws.Range("H34").Value = 4333191
$ws.Range("I34").Value = 1627.3846
$ws.Range("J34").Value = 8778743
$ws.Range("K34").Value = 1627.3846
$ws.Range("L34").Value = 8778743
$ws.Range("M34").Value = -1425.3846
$ws.Range("N34").Value = -8779147

$ws.Range("H35").Value = 1288.6666
$ws.Range("I35").Value = 1288.6666
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1288.6666
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -994.6666
$ws.Range("N35").ClearContents()

$ws.Range("H56").Value = 9000
$ws.Range("I56").Value = 9000
$ws.Range("K56").Value = 9000
$ws.Range("M56").Value = -8155

$ws.Range("H99").Value = 2615.6333
$ws.Range("I99").Value = 2169
$ws.Range("J99").Value = 2874.2104
$ws.Range("K99").Value = 2169
$ws.Range("L99").Value = 2874.2104
$ws.Range("M99").Value = -671
$ws.Range("N99").Value = -5870.2104

$ws.Range("H105").Value = 2933.1904
$ws.Range("I105").Value = 2906.4666
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2906.4666
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1159.4666
$ws.Range("N105").Value = -6494

$ws.Range("H113").Value = 929.8889
$ws.Range("I113").Value = 920.3125
$ws.Range("J113").Value = 1006.5
$ws.Range("K113").Value = 920.3125
$ws.Range("L113").Value = 1006.5
$ws.Range("M113").Value = 1249.6875
$ws.Range("N113").Value = -5346.5

$ws.Range("H126").Value = 2615.6333
$ws.Range("I126").Value = 2169
$ws.Range("J126").Value = 2874.2104
$ws.Range("K126").Value = 6507
$ws.Range("L126").Value = 8622.6312
$ws.Range("M126").Value = -4037
$ws.Range("N126").Value = -13562.6312

$ws.Range("H132").Value = 32409.63
$ws.Range("I132").Value = 1101.0968
$ws.Range("J132").Value = 97113.92999999999
$ws.Range("K132").Value = 3303.2904
$ws.Range("L132").Value = 291341.79
$ws.Range("M132").Value = -773.2903999999999
$ws.Range("N132").Value = -296401.79

$ws.Range("H134").Value = 431578.47
$ws.Range("I134").Value = 486997.53
$ws.Range("J134").Value = 201985.28
$ws.Range("K134").Value = 1460992.59
$ws.Range("L134").Value = 605955.84
$ws.Range("M134").Value = -1458457.59
$ws.Range("N134").Value = -611025.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4360.3706
$ws.Range("I113").Value = 8151
$ws.Range("J113").Value = 840.5
$ws.Range("K113").Value = 24453
$ws.Range("L113").Value = 2521.5
$ws.Range("M113").Value = -22283
$ws.Range("N113").Value = -6861.5

$ws.Range("H121").Value = 86765.5
$ws.Range("I121").Value = 191.66667
$ws.Range("J121").Value = 130052.414
$ws.Range("K121").Value = 575.00001
$ws.Range("L121").Value = 390157.242
$ws.Range("M121").Value = 734.99999
$ws.Range("N121").Value = -392777.242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4796.278
$ws.Range("I70").Value = 4855.567
$ws.Range("J70").Value = 4499.8335
$ws.Range("K70").Value = 4855.567
$ws.Range("L70").Value = 4499.8335
$ws.Range("M70").Value = -4585.567
$ws.Range("N70").Value = -5039.8335

$ws.Range("H73").Value = 4796.278
$ws.Range("I73").Value = 4855.567
$ws.Range("J73").Value = 4499.8335
$ws.Range("K73").Value = 4855.567
$ws.Range("L73").Value = 4499.8335
$ws.Range("M73").Value = -3919.567
$ws.Range("N73").Value = -6371.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 753850.0600000001
$ws.Range("I136").Value = 1061476.6
$ws.Range("K136").Value = 3184429.8
$ws.Range("M136").Value = -3181879.8
